$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge row 143 ('20180423_01_093_094') into row 144 ('20180423_01_095_100'):
# update row 144's value to the merged label, then delete row 143.
$ws.Range("A144").Value = "20180423_01_093_100"
$ws.Rows("143").Delete()

# Update the selection to match the target state.
$ws.Range("B142").Select()
